# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for rows 2-53 on the active
# worksheet to the freshly (re)computed strike-count (K) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, in row order starting at row 2 through row 53.
$kValues = @(
    0,0,2,1,2,0,1,3,2,1,
    0,0,4,0,1,1,1,0,2,1,
    1,4,1,1,1,0,1,1,1,2,
    3,3,3,1,1,1,1,0,1,1,
    0,3,2,0,1,0,2,1,2,0,
    1,1
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
